$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking strings
# like "1.002" or "0.5070" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.931.96"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.637.30"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "214.29"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "0.5070"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "0.2574"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Value = "0.06361"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "19.76"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D11").Value = "0.07743"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "4.294"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "1.646.81"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "0.5456"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "0.0₅7734"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "64.14"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "25.945.27"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "4.452"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "196.10"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "9.934"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "6.141"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "1.890"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "142.98"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").Value = "0.1248"
$ws.Range("E26").Value = "  +9.25%  "
$ws.Range("D27").Value = "6.844"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").Value = "15.60"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").Value = "1.237"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "0.04883"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("D31").Value = "3.247"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "3.202"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "1.552"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").Value = "2.369"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "0.9132"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("D36").Value = "2.570"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").Value = "0.5516"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").Value = "1.123.97"
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D39").Value = "0.01566"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "5.604"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "0.8039"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").Value = "98.49"
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("D44").Value = "0.0₈121"
$ws.Range("E44").Value = "  -9.61%  "
$ws.Range("D45").Value = "1.770.89"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "0.4489"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "55.02"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").Value = "0.05181"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").Value = "7.521"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -0.50%  "
